# Updated symbol list with new cryptocurrency price data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.69%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.96"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.51%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.559"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08046"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.85%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.897"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.53%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.278"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.80%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9451"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.06%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.540"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-10.76%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1168"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.46%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.98%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09682"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.90%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04366"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.59%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1067"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.22%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001277"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.48%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005970"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.74%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.398"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.11%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.07"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "15.21%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1379"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.64%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.54%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04199"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-4.40%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.47%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004287"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.63%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001261"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.14%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003996"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.23%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02647"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-4.40%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05517"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.75%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007578"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-4.66%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1394"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.59%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007982"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-18.33%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002010"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.57%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008836"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-8.53%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006901"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-6.16%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.30%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.002274"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.23%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.007424"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "114.62%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.30%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.30%"
